$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 new values
$ws.Range("A14").Value = 111798760
$ws.Range("B14").Value = 90709
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 5448
$ws.Range("F14").Value = "Svartvit taggsvamp"
$ws.Range("G14").Value = "Phellodon connatus"
$ws.Range("H14").Value = "(Schultz) nom.prov"
$ws.Range("AF14").Value = ""
$ws.Range("AI14").Value = ""
$ws.Range("AR14").Value = ""

# Row 16 new values
$ws.Range("A16").Value = 111798795
$ws.Range("B16").Value = 81076
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 5046
$ws.Range("F16").Value = "Grön jordtunga"
$ws.Range("G16").Value = "Microglossum viride"
$ws.Range("H16").Value = "(Pers.:Fr.) Gillet"
$ws.Range("AF16").Value = ""
